$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 454, shifting existing rows 454:539 down to 455:540
$ws.Rows.Item(454).Insert()

# Copy the date cell number format from the row below (now row 455) so the
# newly inserted row's date cell keeps the same date style.
$ws.Cells.Item(455, 4).Copy()
$ws.Cells.Item(454, 4).PasteSpecial(-4122)  # xlPasteFormats

# Fill the new row 454 with the weekly record's values.
$ws.Cells.Item(454, 1).Value = 4                                          # A - Mercado ID
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"          # B - Mercado
$ws.Cells.Item(454, 3).Value = "Los Lagos"                                # C - Region
$ws.Cells.Item(454, 4).Value = 45258                                      # D - Fecha
$ws.Cells.Item(454, 5).Value = 10                                         # E - Codreg
$ws.Cells.Item(454, 6).Value = 100112043                                  # F - Categoria ID
$ws.Cells.Item(454, 7).Value = "Pepino ensalada"                          # G - Categoria
$ws.Cells.Item(454, 8).Value = "Sin especificar"                          # H - Variedad
$ws.Cells.Item(454, 9).Value = "Primera"                                  # I - Calidad
$ws.Cells.Item(454, 10).Value = 600                                       # J - Volumen
$ws.Cells.Item(454, 11).Value = 19500                                     # K - Precio minimo
$ws.Cells.Item(454, 12).Value = 21000                                     # L - Precio maximo
$ws.Cells.Item(454, 13).Value = 20250                                     # M - Precio promedio ponderado
$ws.Cells.Item(454, 14).Value = "$/caja 60 unidades"                      # N - Unidad de comercializacion
$ws.Cells.Item(454, 15).Value = "Región de Arica y Parinacota"            # O - Origen
$ws.Cells.Item(454, 16).Value = 338                                       # P - Precio $/Kg
$ws.Cells.Item(454, 17).Value = 60                                        # Q - Kg o Unidades
$ws.Cells.Item(454, 18).Value = "Hortaliza"                               # R - Clasificacion
